$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D ("Runmode") was switched from "N" to "Y" for every data row
# (row 59 was already "Y"); this also frees up the now-unused "N" shared
# string, which disappears from the workbook on save.
for ($r = 2; $r -le 63; $r++) {
    $ws.Range("D$r").Value = "Y"
}

# Selection moved back to D2 and the stored scroll position (topLeftCell)
# was reset to the top of the sheet.
$ws.Range("D2").Select() | Out-Null
